$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.027335366613918
$ws.Range("D2").Value = 1.034881058108953
$ws.Range("E2").Value = 1.027387393856686
$ws.Range("F2").Value = 1.042010425207703
$ws.Range("J2").Value = 1.032493811895248
$ws.Range("K2").Value = 1.037679260548263
$ws.Range("L2").Value = 1.030207270465785
$ws.Range("M2").Value = 1.044788321311974
$ws.Range("N2").Value = 1.033960071536853
$ws.Range("C3").Value = 1.028727479668757
$ws.Range("D3").Value = 1.036174129765366
$ws.Range("E3").Value = 1.028585289557304
$ws.Range("F3").Value = 1.043473244770951
$ws.Range("J3").Value = 1.033523793130702
$ws.Range("K3").Value = 1.038780324607271
$ws.Range("L3").Value = 1.031211825966323
$ws.Range("M3").Value = 1.046060176221303
$ws.Range("N3").Value = 1.034991515463784
$ws.Range("C4").Value = 1.029627962063278
$ws.Range("D4").Value = 1.037010819290605
$ws.Range("E4").Value = 1.029360361741199
$ws.Range("F4").Value = 1.044420050316334
$ws.Range("J4").Value = 1.034189557042109
$ws.Range("K4").Value = 1.039492230414935
$ws.Range("L4").Value = 1.031861234186113
$ws.Range("M4").Value = 1.046882906415015
$ws.Range("N4").Value = 1.035658224836309
$ws.Range("C5").Value = 1.030006456032151
$ws.Range("D5").Value = 1.037362563931431
$ws.Range("E5").Value = 1.029686194208471
$ws.Range("F5").Value = 1.044818154975199
$ws.Range("J5").Value = 1.034469279179316
$ws.Range("K5").Value = 1.039791385866541
$ws.Range("L5").Value = 1.032134103209044
$ws.Range("M5").Value = 1.04722872753168
$ws.Range("N5").Value = 1.035938344211035
$ws.Range("C6").Value = 1.030070002940091
$ws.Range("D6").Value = 1.037421623535709
$ws.Range("E6").Value = 1.029740902545407
$ws.Range("F6").Value = 1.044885002595281
$ws.Range("J6").Value = 1.034516236174855
$ws.Range("K6").Value = 1.039841607873909
$ws.Range("L6").Value = 1.032179910866666
$ws.Range("M6").Value = 1.047286789329636
$ws.Range("N6").Value = 1.035985367890894
$ws.Range("C7").Value = 1.029633019782616
$ws.Range("D7").Value = 1.037015519313573
$ws.Range("E7").Value = 1.029364715557515
$ws.Range("F7").Value = 1.044425369538983
$ws.Range("J7").Value = 1.034193295347355
$ws.Range("K7").Value = 1.039496228250307
$ws.Range("L7").Value = 1.031864880830391
$ws.Range("M7").Value = 1.046887527503509
$ws.Range("N7").Value = 1.035661968450378
$ws.Range("C8").Value = 1.027805904410676
$ws.Range("D8").Value = 1.035318062972578
$ws.Range("E8").Value = 1.027792239882794
$ws.Range("F8").Value = 1.042504740103737
$ws.Range("J8").Value = 1.032842046354945
$ws.Range("K8").Value = 1.038051487541046
$ws.Range("L8").Value = 1.030546892041971
$ws.Range("M8").Value = 1.045218203762037
$ws.Range("N8").Value = 1.034308800529408
$ws.Range("C9").Value = 1.024583718032479
$ws.Range("D9").Value = 1.032326648677793
$ws.Range("E9").Value = 1.025020835944709
$ws.Range("F9").Value = 1.039122160523541
$ws.Range("J9").Value = 1.03045544471044
$ws.Range("K9").Value = 1.035501256147781
$ws.Range("L9").Value = 1.028219650109889
$ws.Range("M9").Value = 1.042274586937265
$ws.Range("N9").Value = 1.031918809636821
$ws.Range("C10").Value = 1.022433556368671
$ws.Range("D10").Value = 1.030331945821978
$ws.Range("E10").Value = 1.023172691857871
$ws.Range("F10").Value = 1.036868042639147
$ws.Range("J10").Value = 1.028860455335134
$ws.Range("K10").Value = 1.033797918677249
$ws.Range("L10").Value = 1.026664755764723
$ws.Range("M10").Value = 1.040310542663704
$ws.Range("N10").Value = 1.030321555193654
$ws.Range("C11").Value = 1.021501961062661
$ws.Range("D11").Value = 1.029468062673377
$ws.Range("E11").Value = 1.022372250587148
$ws.Range("F11").Value = 1.035892142708853
$ws.Range("J11").Value = 1.028168831217057
$ws.Range("K11").Value = 1.033059551325
$ws.Range("L11").Value = 1.025990620762613
$ws.Range("M11").Value = 1.039459646155184
$ws.Range("N11").Value = 1.029628948889998
$ws.Range("C12").Value = 1.021155834842575
$ws.Range("D12").Value = 1.029147148739438
$ws.Range("E12").Value = 1.022074899577407
$ws.Range("F12").Value = 1.035529666700691
$ws.Range("J12").Value = 1.027911779487191
$ws.Range("K12").Value = 1.032785162737395
$ws.Range("L12").Value = 1.025740084775816
$ws.Range("M12").Value = 1.039143512426874
$ws.Range("N12").Value = 1.029371532117192
$ws.Range("C13").Value = 1.021230084271147
$ws.Range("D13").Value = 1.029215987267473
$ws.Range("E13").Value = 1.022138683868013
$ws.Range("F13").Value = 1.035607418369639
$ws.Range("J13").Value = 1.027966924887085
$ws.Range("K13").Value = 1.032844025794328
$ws.Range("L13").Value = 1.025793831638828
$ws.Range("M13").Value = 1.039211327561785
$ws.Range("N13").Value = 1.029426755829879
$ws.Range("C14").Value = 1.021473352036747
$ws.Range("D14").Value = 1.029441536451704
$ws.Range("E14").Value = 1.022347672125561
$ws.Range("F14").Value = 1.035862180020939
$ws.Range("J14").Value = 1.028147586359625
$ws.Range("K14").Value = 1.033036872875255
$ws.Range("L14").Value = 1.025969914098747
$ws.Range("M14").Value = 1.03943351595779
$ws.Range("N14").Value = 1.029607673862432
$ws.Range("C15").Value = 1.021623225199579
$ws.Range("D15").Value = 1.029580500729455
$ws.Range("E15").Value = 1.022476432395782
$ws.Range("F15").Value = 1.036019149125892
$ws.Range("J15").Value = 1.028258877628021
$ws.Range("K15").Value = 1.033155675493107
$ws.Range("L15").Value = 1.026078386680774
$ws.Range("M15").Value = 1.039570403757194
$ws.Range("N15").Value = 1.029719123177195
$ws.Range("C16").Value = 1.02249537095221
$ws.Range("D16").Value = 1.030389275054379
$ws.Range("E16").Value = 1.023225810315806
$ws.Range("F16").Value = 1.036932812547738
$ws.Range("J16").Value = 1.028906335092201
$ws.Range("K16").Value = 1.033846904238985
$ws.Range("L16").Value = 1.026709477544683
$ws.Range("M16").Value = 1.040367003912549
$ws.Range("N16").Value = 1.030367500105239
$ws.Range("C17").Value = 1.023042290525656
$ws.Range("D17").Value = 1.030896550708069
$ws.Range("E17").Value = 1.023695823532942
$ws.Range("F17").Value = 1.0375059651788
$ws.Range("J17").Value = 1.029312201989157
$ws.Range("K17").Value = 1.034280273659996
$ws.Range("L17").Value = 1.027105112521694
$ws.Range("M17").Value = 1.040866565566367
$ws.Range("N17").Value = 1.030773943379742
$ws.Range("C18").Value = 1.023361245723391
$ws.Range("D18").Value = 1.031192420542775
$ws.Range("E18").Value = 1.023969956938422
$ws.Range("F18").Value = 1.037840290067524
$ws.Range("J18").Value = 1.029548842692931
$ws.Range("K18").Value = 1.034532972750476
$ws.Range("L18").Value = 1.027335797419521
$ws.Range("M18").Value = 1.041157908113015
$ws.Range("N18").Value = 1.031010920140458
$ws.Range("C19").Value = 1.023469992267761
$ws.Range("D19").Value = 1.031293302027151
$ws.Range("E19").Value = 1.024063426523683
$ws.Range("F19").Value = 1.03795428886812
$ws.Range("J19").Value = 1.029629515134556
$ws.Range("K19").Value = 1.034619123514041
$ws.Range("L19").Value = 1.0274144411618
$ws.Range("M19").Value = 1.041257241193648
$ws.Range("N19").Value = 1.031091707146203
$ws.Range("C20").Value = 1.022983616834215
$ws.Range("D20").Value = 1.030842126466818
$ws.Range("E20").Value = 1.023645397388189
$ws.Range("F20").Value = 1.037444469851583
$ws.Range("J20").Value = 1.029268666126463
$ws.Range("K20").Value = 1.034233785318249
$ws.Range("L20").Value = 1.02706267317257
$ws.Range("M20").Value = 1.040812971845881
$ws.Range("N20").Value = 1.030730345691129
$ws.Range("C21").Value = 1.021401718282093
$ws.Range("D21").Value = 1.029375118684259
$ws.Range("E21").Value = 1.022286131194175
$ws.Range("F21").Value = 1.03578715864874
$ws.Range("J21").Value = 1.028094390279943
$ws.Range("K21").Value = 1.032980087730962
$ws.Range("L21").Value = 1.025918065883094
$ws.Range("M21").Value = 1.039368089064449
$ws.Range("N21").Value = 1.029554402238215
$ws.Range("C22").Value = 1.020406590792213
$ws.Range("D22").Value = 1.028452580996317
$ws.Range("E22").Value = 1.021431322430793
$ws.Range("F22").Value = 1.034745235042275
$ws.Range("J22").Value = 1.027355196534637
$ws.Range("K22").Value = 1.03219110674453
$ws.Range("L22").Value = 1.025197638956193
$ws.Range("M22").Value = 1.038459209809446
$ws.Range("N22").Value = 1.028814158753005
$ws.Range("C23").Value = 1.020934178464496
$ws.Range("D23").Value = 1.028941653519685
$ws.Range("E23").Value = 1.021884491537471
$ws.Range("F23").Value = 1.035297571435952
$ws.Range("J23").Value = 1.027747141878149
$ws.Range("K23").Value = 1.032609431254528
$ws.Range("L23").Value = 1.025579624834042
$ws.Range("M23").Value = 1.038941065744513
$ws.Range("N23").Value = 1.02920666070386
$ws.Range("C24").Value = 1.023010129115449
$ws.Range("D24").Value = 1.030866718486321
$ws.Range("E24").Value = 1.023668182845848
$ws.Range("F24").Value = 1.037472256897304
$ws.Range("J24").Value = 1.02928833840087
$ws.Range("K24").Value = 1.034254791639458
$ws.Range("L24").Value = 1.027081849941721
$ws.Range("M24").Value = 1.040837188674796
$ws.Range("N24").Value = 1.030750045902421
$ws.Range("C25").Value = 1.025417067088078
$ws.Range("D25").Value = 1.033100058618121
$ws.Range("E25").Value = 1.025737390675231
$ws.Range("F25").Value = 1.039996452642755
$ws.Range("J25").Value = 1.031073114602026
$ws.Range("K25").Value = 1.036161097295538
$ws.Range("L25").Value = 1.028821883351618
$ws.Range("M25").Value = 1.04303585180819
$ws.Range("N25").Value = 1.032537356690501

Write-Output "Applied 216 cell updates"